$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.57915860414505
$ws.Range("B1").Value = 1.059358835220337
$ws.Range("C1").Value = 3.79627537727356
$ws.Range("D1").Value = 3.061160564422607
$ws.Range("E1").Value = 0.8158451318740845
